$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.970.95'
$ws.Range('E2').Value = '  -1.49%  '

# Row 3
$ws.Range('D3').Value = '3.383.59'
$ws.Range('E3').Value = '  -0.91%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.92%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.42%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
$ws.Range('D8').Value = '3.387.67'
$ws.Range('E8').Value = '  -0.72%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.471'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.35%  '

# Row 10
$ws.Range('E10').Value = '  +0.41%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.124'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.30%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.388'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.14%  '

# Row 13
$ws.Range('D13').Value = '3.954.21'
$ws.Range('E13').Value = '  -1.10%  '

# Row 14
$ws.Range('E14').Value = '  +1.12%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.45'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.81%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000173'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.61%  '

# Row 17
$ws.Range('D17').Value = '3.377.89'
$ws.Range('E17').Value = '  -1.03%  '

# Row 18
$ws.Range('D18').Value = '61.080.61'
$ws.Range('E18').Value = '  -1.39%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.02'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.70%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.84'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.20%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.38'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.42%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '377.66'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.38%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.557'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.81%  '

# Row 24
$ws.Range('D24').Value = '3.506.10'
$ws.Range('E24').Value = '  -1.26%  '

# Row 25
$ws.Range('E25').Value = '  +0.04%  '

# Row 26
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000125'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.33%  '

# Row 27
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '71.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.69%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +11.06%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.55%  '

# Row 30
$ws.Range('E30').Value = '  +4.37%  '

# Row 31
$ws.Range('E31').Value = '  -0.17%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.16'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.87%  '

# Row 33
$ws.Range('E33').Value = '  -1.03%  '

# Row 34
$ws.Range('E34').Value = '  -0.04%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.79'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.78%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.23'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.00%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.64%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.99%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '164.81'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.31%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0756'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.73%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.12%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.772'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.80%  '

# Row 43
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.70'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.65%  '

# Row 44
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.81%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.71%  '

# Row 46
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '41.52'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.31%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '23.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.14%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.40'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.16%  '

# Row 49
$ws.Range('E49').Value = '  -2.61%  '

# Row 50
$ws.Range('D50').Value = '2.373.29'
$ws.Range('E50').Value = '  -1.05%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.42'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.81%  '
